$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.317.34"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.34%  "

$ws.Range("D3").Value = "'1.840.51"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.73%  "

$ws.Range("D4").Value = "'0.9990"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "'239.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.54%  "

$ws.Range("D6").Value = "'0.6264"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.18%  "

$ws.Range("D7").Value = "'1.000"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").Value = "'0.07410"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.38%  "

$ws.Range("D9").Value = "'0.2890"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.44%  "

$ws.Range("D10").Value = "'24.71"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.65%  "

$ws.Range("D11").Value = "'0.07733"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.27%  "

$ws.Range("D12").Value = "'1.831.12"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.23%  "

$ws.Range("D13").Value = "'4.973"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.13%  "

$ws.Range("D14").Value = "'0.6759"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.99%  "

$ws.Range("D15").Value = "'0.00001022"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.70%  "

$ws.Range("D16").Value = "'81.95"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.70%  "

$ws.Range("D17").Value = "'6.234"
$ws.Range("D17").Style = "Normal"

$ws.Range("D18").Value = "'29.305.98"
$ws.Range("D18").Style = "Normal"

$ws.Range("D19").Value = "'228.54"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.86%  "

$ws.Range("D20").Value = "'12.27"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.06%  "

$ws.Range("D21").Value = "'0.9999"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.03%  "

$ws.Range("D22").Value = "'7.404"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.15%  "

$ws.Range("D23").Value = "'1.000"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.04%  "

$ws.Range("D24").Value = "'158.83"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.38%  "

$ws.Range("D25").Value = "'8.455"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.12%  "

$ws.Range("D26").Value = "'0.1348"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.67%  "

$ws.Range("E27").Value = "  -1.66%  "

$ws.Range("D28").Value = "'0.06632"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +16.35%  "

$ws.Range("E29").Value = "  +2.44%  "

$ws.Range("D30").Value = "'1.485"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.44%  "

$ws.Range("D31").Value = "'4.057"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.88%  "

$ws.Range("E32").Value = "  +0.10%  "

$ws.Range("E33").Value = "  +0.05%  "

$ws.Range("D34").Value = "'1.135"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.98%  "

$ws.Range("D35").Value = "'0.6908"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.27%  "

$ws.Range("D36").Value = "'2.567"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.54%  "

$ws.Range("D37").Value = "'0.01851"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.12%  "

$ws.Range("D38").Value = "'2.822"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.52%  "

$ws.Range("D39").Value = "'1.242.57"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.46%  "

$ws.Range("D40").Value = "'6.747"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.28%  "

$ws.Range("D41").Value = "'0.9317"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.08%  "

$ws.Range("D42").Value = "'0.9996"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.06%  "

$ws.Range("D43").Value = "'1.980.82"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.58%  "

$ws.Range("D44").Value = "'100.57"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.42%  "

$ws.Range("D45").Value = "'65.43"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.93%  "

$ws.Range("D46").Value = "'7.025"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.86%  "

$ws.Range("D47").Value = "'1.706"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.64%  "

$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").Value = "'0.1148"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.96%  "

$ws.Range("D49").Value = "'8.960"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.70%  "

$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "'0.00000000115"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.25%  "

$ws.Range("D51").Value = "'0.3881"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.31%  "
